$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A1").Value = "venkat1"
$ws.Range("A2").Value = "venkat2"
$ws.Range("A3").Value = "venkat3"
$ws.Range("A4").Value = "venkat4"

$ws.Range("A5").Select()
